$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" note in A2
$ws.Range("A2").Value = "last updated: 02.07.2025"

# Insert two new rows at 25:26 for the new GC autosampler / guard column fields
# (shifts the former rows 26-41 down to 28-43)
$ws.Rows("25:26").Insert()

# The insert carries formatting from the row above into stray cells in the
# newly inserted rows (e.g. B/D) - clear those so only the intended cells
# hold content, matching the template's sparse layout.
$ws.Range("A25:H26").Clear()

# Populate the new "guard column" and "autosampler" rows (guard column
# entered first, matching the authored shared-string ordering)
$ws.Range("F26").Value = "guard column"
$ws.Range("F25").Value = "autopsampler model"
$ws.Range("G26").Value = "chromatography guard column model"
$ws.Range("G25").Value = "chromatography autosampler model"

# Restore the selection to match the authored workbook
$ws.Range("G31").Select() | Out-Null
